$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = '71.489.94'
$ws.Cells.Item(2, 5).Value = '  -1.64%  '

# Row 3
$ws.Cells.Item(3, 4).Value = '3.971.73'
$ws.Cells.Item(3, 5).Value = '  -2.10%  '

# Row 4
$ws.Cells.Item(4, 5).Value = '  -0.05%  '

# Row 5
$ws.Cells.Item(5, 4).NumberFormat = '@'
$ws.Cells.Item(5, 4).Value = '543.11'
$ws.Cells.Item(5, 5).Value = '  +4.64%  '

# Row 6
$ws.Cells.Item(6, 4).NumberFormat = '@'
$ws.Cells.Item(6, 4).Value = '149.30'
$ws.Cells.Item(6, 5).Value = '  +1.51%  '

# Row 7
$ws.Cells.Item(7, 4).Value = '3.967.18'
$ws.Cells.Item(7, 5).Value = '  -2.03%  '

# Row 8
$ws.Cells.Item(8, 5).Value = '  -6.82%  '

# Row 9
$ws.Cells.Item(9, 5).Value = '  +0.05%  '

# Row 10
$ws.Cells.Item(10, 4).NumberFormat = '@'
$ws.Cells.Item(10, 4).Value = '0.742'
$ws.Cells.Item(10, 5).Value = '  -4.03%  '

# Row 11
$ws.Cells.Item(11, 5).Value = '  -5.55%  '

# Row 12
$ws.Cells.Item(12, 4).NumberFormat = '@'
$ws.Cells.Item(12, 4).Value = '56.98'
$ws.Cells.Item(12, 5).Value = '  +19.39%  '

# Row 13
$ws.Cells.Item(13, 4).NumberFormat = '@'
$ws.Cells.Item(13, 4).Value = '0.0000319'
$ws.Cells.Item(13, 5).Value = '  -2.47%  '

# Row 14
$ws.Cells.Item(14, 4).NumberFormat = '@'
$ws.Cells.Item(14, 4).Value = '10.73'
$ws.Cells.Item(14, 5).Value = '  -4.11%  '

# Row 15
$ws.Cells.Item(15, 4).Value = '4.603.04'
$ws.Cells.Item(15, 5).Value = '  -2.22%  '

# Row 16
$ws.Cells.Item(16, 4).Value = '3.972.72'
$ws.Cells.Item(16, 5).Value = '  -2.17%  '

# Row 17
$ws.Cells.Item(17, 4).NumberFormat = '@'
$ws.Cells.Item(17, 4).Value = '13.92'
$ws.Cells.Item(17, 5).Value = '  -1.29%  '

# Row 18
$ws.Cells.Item(18, 4).NumberFormat = '@'
$ws.Cells.Item(18, 4).Value = '20.53'
$ws.Cells.Item(18, 5).Value = '  -3.32%  '

# Row 19
$ws.Cells.Item(19, 5).Value = '  -1.21%  '

# Row 20
$ws.Cells.Item(20, 4).NumberFormat = '@'
$ws.Cells.Item(20, 4).Value = '1.18'
$ws.Cells.Item(20, 5).Value = '  -2.83%  '

# Row 21
$ws.Cells.Item(21, 4).Value = '71.389.84'
$ws.Cells.Item(21, 5).Value = '  -1.65%  '

# Row 22
$ws.Cells.Item(22, 4).NumberFormat = '@'
$ws.Cells.Item(22, 4).Value = '427.05'
$ws.Cells.Item(22, 5).Value = '  -3.65%  '

# Row 23
$ws.Cells.Item(23, 2).Value = 'ImmutableX'
$ws.Cells.Item(23, 3).Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Cells.Item(23, 4).NumberFormat = '@'
$ws.Cells.Item(23, 4).Value = '3.59'
$ws.Cells.Item(23, 5).Value = '  +0.58%  '

# Row 24
$ws.Cells.Item(24, 2).Value = 'Litecoin'
$ws.Cells.Item(24, 3).Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Cells.Item(24, 4).NumberFormat = '@'
$ws.Cells.Item(24, 4).Value = '97.54'
$ws.Cells.Item(24, 5).Value = '  -7.00%  '

# Row 25
$ws.Cells.Item(25, 4).NumberFormat = '@'
$ws.Cells.Item(25, 4).Value = '4.22'
$ws.Cells.Item(25, 5).Value = '  +5.38%  '

# Row 26
$ws.Cells.Item(26, 4).NumberFormat = '@'
$ws.Cells.Item(26, 4).Value = '14.45'
$ws.Cells.Item(26, 5).Value = '  -2.32%  '

# Row 27
$ws.Cells.Item(27, 4).NumberFormat = '@'
$ws.Cells.Item(27, 4).Value = '11.52'
$ws.Cells.Item(27, 5).Value = '  +0.75%  '

# Row 28
$ws.Cells.Item(28, 4).NumberFormat = '@'
$ws.Cells.Item(28, 4).Value = '10.77'
$ws.Cells.Item(28, 5).Value = '  -2.22%  '

# Row 29
$ws.Cells.Item(29, 4).NumberFormat = '@'
$ws.Cells.Item(29, 4).Value = '3.76'
$ws.Cells.Item(29, 5).Value = '  +14.76%  '

# Row 30
$ws.Cells.Item(30, 4).NumberFormat = '@'
$ws.Cells.Item(30, 4).Value = '5.91'
$ws.Cells.Item(30, 5).Value = '  +1.70%  '

# Row 31
$ws.Cells.Item(31, 4).NumberFormat = '@'
$ws.Cells.Item(31, 4).Value = '36.65'
$ws.Cells.Item(31, 5).Value = '  -2.99%  '

# Row 32
$ws.Cells.Item(32, 4).NumberFormat = '@'
$ws.Cells.Item(32, 4).Value = '7.80'
$ws.Cells.Item(32, 5).Value = '  +14.26%  '

# Row 33
$ws.Cells.Item(33, 4).NumberFormat = '@'
$ws.Cells.Item(33, 4).Value = '51.66'
$ws.Cells.Item(33, 5).Value = '  +21.10%  '

# Row 34
$ws.Cells.Item(34, 4).NumberFormat = '@'
$ws.Cells.Item(34, 4).Value = '695.86'
$ws.Cells.Item(34, 5).Value = '  +2.21%  '

# Row 35
$ws.Cells.Item(35, 2).Value = 'Hedera'
$ws.Cells.Item(35, 3).Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Cells.Item(35, 4).NumberFormat = '@'
$ws.Cells.Item(35, 4).Value = '0.131'
$ws.Cells.Item(35, 5).Value = '  +0.82%  '

# Row 36
$ws.Cells.Item(36, 2).Value = 'Cosmos'
$ws.Cells.Item(36, 3).Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Cells.Item(36, 4).NumberFormat = '@'
$ws.Cells.Item(36, 4).Value = '13.42'
$ws.Cells.Item(36, 5).Value = '  -1.80%  '

# Row 37
$ws.Cells.Item(37, 4).NumberFormat = '@'
$ws.Cells.Item(37, 4).Value = '65.05'
$ws.Cells.Item(37, 5).Value = '  -2.94%  '

# Row 38
$ws.Cells.Item(38, 4).NumberFormat = '@'
$ws.Cells.Item(38, 4).Value = '0.438'
$ws.Cells.Item(38, 5).Value = '  +2.24%  '

# Row 39
$ws.Cells.Item(39, 2).Value = 'PEPE'
$ws.Cells.Item(39, 3).Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Cells.Item(39, 4).Value = '0.0₃0832'
$ws.Cells.Item(39, 5).Value = '  -3.80%  '

# Row 40
$ws.Cells.Item(40, 2).Value = 'Kaspa'
$ws.Cells.Item(40, 3).Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Cells.Item(40, 4).NumberFormat = '@'
$ws.Cells.Item(40, 4).Value = '0.151'
$ws.Cells.Item(40, 5).Value = '  +0.15%  '

# Row 41
$ws.Cells.Item(41, 5).Value = '  -3.06%  '

# Row 42
$ws.Cells.Item(42, 5).Value = '  +0.23%  '

# Row 43
$ws.Cells.Item(43, 5).Value = '  +0.17%  '

# Row 44
$ws.Cells.Item(44, 4).NumberFormat = '@'
$ws.Cells.Item(44, 4).Value = '3.27'
$ws.Cells.Item(44, 5).Value = '  +0.80%  '

# Row 45
$ws.Cells.Item(45, 4).NumberFormat = '@'
$ws.Cells.Item(45, 4).Value = '0.0485'
$ws.Cells.Item(45, 5).Value = '  -2.48%  '

# Row 46
$ws.Cells.Item(46, 5).Value = '  -6.33%  '

# Row 47
$ws.Cells.Item(47, 4).NumberFormat = '@'
$ws.Cells.Item(47, 4).Value = '2.71'
$ws.Cells.Item(47, 5).Value = '  +0.49%  '

# Row 48
$ws.Cells.Item(48, 4).NumberFormat = '@'
$ws.Cells.Item(48, 4).Value = '9.80'
$ws.Cells.Item(48, 5).Value = '  +6.95%  '

# Row 49
$ws.Cells.Item(49, 4).NumberFormat = '@'
$ws.Cells.Item(49, 4).Value = '3.37'
$ws.Cells.Item(49, 5).Value = '  -3.84%  '

# Row 50
$ws.Cells.Item(50, 5).Value = '  -1.94%  '

# Row 51
$ws.Cells.Item(51, 5).Value = '  +2.42%  '
